# Adds a new job posting as the second row of Sheet1 (right after the header
# row), shifting all existing job rows down by one. This mirrors how a new
# job entry was inserted at the top of the jobs list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new blank row at row 2, pushing current rows 2-9 down to 3-10.
$ws.Rows.Item(2).Insert()

# The inserted row inherits the header row's formatting; strip that back off
# so the new row matches the plain (unstyled) data rows, then drop the
# leftover empty cell placeholders.
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(2).ClearContents()

# Populate the new job entry (note: column A / S.No intentionally left blank,
# matching the source data for this row).
$ws.Cells.Item(2, 2).Value  = "Samkwang Mobile"
$ws.Cells.Item(2, 3).Value  = "kasna Gr Noida"
$ws.Cells.Item(2, 4).Value  = "Not in use"
$ws.Cells.Item(2, 5).Value  = "Not in use"
$ws.Cells.Item(2, 6).Value  = "SSR contractor में Boysकी भर्ती है आवशयकता है samkwang company के लिए"
$ws.Cells.Item(2, 7).Value  = "वर्क लोकेशन (फैक्ट्री) -: Samkwang india electronics pvt Ltd.. Company"
$ws.Cells.Item(2, 8).Value  = "Kasna , Greater Noida."
$ws.Cells.Item(2, 9).Value  = "Department: Assembly"
$ws.Cells.Item(2, 10).Value = "Salary = 10994  and  OT= 106 per hrs"
$ws.Cells.Item(2, 11).Value = "योग्यता: 10th, 12 th,ITI"
$ws.Cells.Item(2, 12).Value = "आवश्यक दस्तावेज़: resume 2 photo , photo copy all documents"
$ws.Cells.Item(2, 13).Value = "अतिरिक्त लाभ: Double Overtime, Attend.Award- 700, Lunch/Canteen Free, Bus Free"
$ws.Cells.Item(2, 14).Value = "Only for Boys"
$ws.Cells.Item(2, 15).Value = "नोट- डॉक्यूमेंट के साथ आधार कार्ड और बैंक अकाउंट लाना अनिवार्य है (बैंक की पासबुक या चेक बुक)"
$ws.Cells.Item(2, 16).Value = "न्यूनतम दूरी- 100 km"
$ws.Cells.Item(2, 17).Value = "इंटरव्यू डेट :- *09/04/2025"
$ws.Cells.Item(2, 18).Value = "Morning- 7:00 बजे"
$ws.Cells.Item(2, 19).Value = "SS Research Solution Services"
$ws.Cells.Item(2, 20).Value = "Contact: Abhishek Rawal-9990294311, Vivek-9368772854, Kush-7253035420, 9758570409 Ashish"
$ws.Cells.Item(2, 21).Value = "Address: Plot No.1, Samkwang India Electronic Private limited, 10, Ecotech III, Greater Noida, Uttar Pradesh 203202"
$ws.Cells.Item(2, 22).Value = "Google map: https://maps.app.goo.gl/kBG75cpD1WuKs6Tg9"

# Update the view: scroll the frozen/visible area over a bit and move the
# active selection, matching the saved view state after the edit.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("W2").Select()
